# forecast update last observation 21-1-25
# Refresh the VCI3M forecast overview for CLUSTER_1: shift the 11
# forecast-horizon date headers (row 1) forward by 14 days and write the
# newly recalculated forecast values for every cluster row (rows 2-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: forecast horizon dates (B1:L1)
$ws.Range("B1").Value = 45678
$ws.Range("C1").Value = 45685
$ws.Range("D1").Value = 45692
$ws.Range("E1").Value = 45699
$ws.Range("F1").Value = 45706
$ws.Range("G1").Value = 45713
$ws.Range("H1").Value = 45720
$ws.Range("I1").Value = 45727
$ws.Range("J1").Value = 45734
$ws.Range("K1").Value = 45741
$ws.Range("L1").Value = 45748

# Row 2: _Average
$ws.Range("B2").Value = 61.9
$ws.Range("C2").Value = 60.7
$ws.Range("D2").Value = 59.5
$ws.Range("E2").Value = 58.5
$ws.Range("F2").Value = 57.5
$ws.Range("G2").Value = 56.8
$ws.Range("H2").Value = 56.3
$ws.Range("I2").Value = 56.1
$ws.Range("J2").Value = 56.2
$ws.Range("K2").Value = 56.6
$ws.Range("L2").Value = 57.2

# Row 3: Abim
$ws.Range("B3").Value = 62.4
$ws.Range("C3").Value = 63.8
$ws.Range("D3").Value = 65
$ws.Range("E3").Value = 66.2
$ws.Range("F3").Value = 67.2
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = 68.5
$ws.Range("I3").Value = 68.8
$ws.Range("J3").Value = 68.8
$ws.Range("K3").Value = 68.5
$ws.Range("L3").Value = 68

# Row 4: Amudat
$ws.Range("B4").Value = 43.4
$ws.Range("C4").Value = 42
$ws.Range("D4").Value = 40.6
$ws.Range("E4").Value = 39.4
$ws.Range("F4").Value = 38.3
$ws.Range("G4").Value = 37.5
$ws.Range("H4").Value = 37
$ws.Range("I4").Value = 36.9
$ws.Range("J4").Value = 37.3
$ws.Range("K4").Value = 38
$ws.Range("L4").Value = 39

# Row 5: Budi
$ws.Range("B5").Value = 65.9
$ws.Range("C5").Value = 64
$ws.Range("D5").Value = 62.2
$ws.Range("E5").Value = 60.6
$ws.Range("F5").Value = 59.2
$ws.Range("G5").Value = 58
$ws.Range("H5").Value = 57.3
$ws.Range("I5").Value = 56.9
$ws.Range("J5").Value = 56.9
$ws.Range("K5").Value = 57.3
$ws.Range("L5").Value = 57.9

# Row 6: Dasenech (Kuraz)
$ws.Range("B6").Value = 62.5
$ws.Range("C6").Value = 64.1
$ws.Range("D6").Value = 65.7
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 68.1
$ws.Range("G6").Value = 68.8
$ws.Range("H6").Value = 69
$ws.Range("I6").Value = 68.8
$ws.Range("J6").Value = 68.3
$ws.Range("K6").Value = 67.4
$ws.Range("L6").Value = 66.3

# Row 7: Gnangatom
$ws.Range("B7").Value = 61
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 59
$ws.Range("E7").Value = 58
$ws.Range("F7").Value = 57.2
$ws.Range("G7").Value = 56.5
$ws.Range("H7").Value = 56.3
$ws.Range("I7").Value = 56.5
$ws.Range("J7").Value = 57.1
$ws.Range("K7").Value = 58.2
$ws.Range("L7").Value = 59.5

# Row 8: Kaabong
$ws.Range("B8").Value = 70.3
$ws.Range("C8").Value = 68.3
$ws.Range("D8").Value = 66.3
$ws.Range("E8").Value = 64.4
$ws.Range("F8").Value = 62.7
$ws.Range("G8").Value = 61.3
$ws.Range("H8").Value = 60.2
$ws.Range("I8").Value = 59.6
$ws.Range("J8").Value = 59.4
$ws.Range("K8").Value = 59.6
$ws.Range("L8").Value = 60.2

# Row 9: Kapoeta East
$ws.Range("B9").Value = 63
$ws.Range("C9").Value = 60.4
$ws.Range("D9").Value = 57.9
$ws.Range("E9").Value = 55.5
$ws.Range("F9").Value = 53.6
$ws.Range("G9").Value = 52.2
$ws.Range("H9").Value = 51.5
$ws.Range("I9").Value = 51.5
$ws.Range("J9").Value = 52.2
$ws.Range("K9").Value = 53.6
$ws.Range("L9").Value = 55.4

# Row 10: Kapoeta North
$ws.Range("B10").Value = 63.7
$ws.Range("C10").Value = 60.2
$ws.Range("D10").Value = 56.8
$ws.Range("E10").Value = 53.7
$ws.Range("F10").Value = 50.9
$ws.Range("G10").Value = 48.7
$ws.Range("H10").Value = 47
$ws.Range("I10").Value = 46
$ws.Range("J10").Value = 45.7
$ws.Range("K10").Value = 46
$ws.Range("L10").Value = 46.8

# Row 11: Kapoeta South
$ws.Range("B11").Value = 74.9
$ws.Range("C11").Value = 72.2
$ws.Range("D11").Value = 69.4
$ws.Range("E11").Value = 66.7
$ws.Range("F11").Value = 64.2
$ws.Range("G11").Value = 62
$ws.Range("H11").Value = 60.3
$ws.Range("I11").Value = 59.1
$ws.Range("J11").Value = 58.5
$ws.Range("K11").Value = 58.5
$ws.Range("L11").Value = 58.9

# Row 12: Karenga
$ws.Range("B12").Value = 66.7
$ws.Range("C12").Value = 65.2
$ws.Range("D12").Value = 63.7
$ws.Range("E12").Value = 62.1
$ws.Range("F12").Value = 60.5
$ws.Range("G12").Value = 58.9
$ws.Range("H12").Value = 57.5
$ws.Range("I12").Value = 56.3
$ws.Range("J12").Value = 55.2
$ws.Range("K12").Value = 54.4
$ws.Range("L12").Value = 53.8

# Row 13: Kotido
$ws.Range("B13").Value = 63.3
$ws.Range("C13").Value = 62
$ws.Range("D13").Value = 60.5
$ws.Range("E13").Value = 59.1
$ws.Range("F13").Value = 57.6
$ws.Range("G13").Value = 56.3
$ws.Range("H13").Value = 55
$ws.Range("I13").Value = 53.9
$ws.Range("J13").Value = 53
$ws.Range("K13").Value = 52.2
$ws.Range("L13").Value = 51.6

# Row 14: Loima
$ws.Range("B14").Value = 58.3
$ws.Range("C14").Value = 58.2
$ws.Range("D14").Value = 58.6
$ws.Range("E14").Value = 59.7
$ws.Range("F14").Value = 61.4
$ws.Range("G14").Value = 63.7
$ws.Range("H14").Value = 66.5
$ws.Range("I14").Value = 69.8
$ws.Range("J14").Value = 73.4
$ws.Range("K14").Value = 77
$ws.Range("L14").Value = 80.6

# Row 15: Moroto
$ws.Range("B15").Value = 64.8
$ws.Range("C15").Value = 61.7
$ws.Range("D15").Value = 58.6
$ws.Range("E15").Value = 55.6
$ws.Range("F15").Value = 52.9
$ws.Range("G15").Value = 50.8
$ws.Range("H15").Value = 49.4
$ws.Range("I15").Value = 48.8
$ws.Range("J15").Value = 49
$ws.Range("K15").Value = 50.1
$ws.Range("L15").Value = 51.8

# Row 16: Nakapiripirit
$ws.Range("B16").Value = 56
$ws.Range("C16").Value = 54.9
$ws.Range("D16").Value = 53.9
$ws.Range("E16").Value = 52.9
$ws.Range("F16").Value = 52.1
$ws.Range("G16").Value = 51.3
$ws.Range("H16").Value = 50.7
$ws.Range("I16").Value = 50.1
$ws.Range("J16").Value = 49.7
$ws.Range("K16").Value = 49.4
$ws.Range("L16").Value = 49.2

# Row 17: Napak
$ws.Range("B17").Value = 63.3
$ws.Range("C17").Value = 62.1
$ws.Range("D17").Value = 60.8
$ws.Range("E17").Value = 59.4
$ws.Range("F17").Value = 58.1
$ws.Range("G17").Value = 56.7
$ws.Range("H17").Value = 55.4
$ws.Range("I17").Value = 54.3
$ws.Range("J17").Value = 53.4
$ws.Range("K17").Value = 52.6
$ws.Range("L17").Value = 52.1

# Row 18: Surma
$ws.Range("B18").Value = 63.6
$ws.Range("C18").Value = 61.3
$ws.Range("D18").Value = 59.2
$ws.Range("E18").Value = 57.3
$ws.Range("F18").Value = 55.8
$ws.Range("G18").Value = 54.7
$ws.Range("H18").Value = 54.2
$ws.Range("I18").Value = 54.3
$ws.Range("J18").Value = 54.9
$ws.Range("K18").Value = 56.1
$ws.Range("L18").Value = 57.6

# Row 19: Turkana
$ws.Range("B19").Value = 55.2
$ws.Range("C19").Value = 56.2
$ws.Range("D19").Value = 57.3
$ws.Range("E19").Value = 58.4
$ws.Range("F19").Value = 59.4
$ws.Range("G19").Value = 60.3
$ws.Range("H19").Value = 61
$ws.Range("I19").Value = 61.5
$ws.Range("J19").Value = 61.8
$ws.Range("K19").Value = 62
$ws.Range("L19").Value = 62

# Row 20: Turkana West
$ws.Range("B20").Value = 64.4
$ws.Range("C20").Value = 64.3
$ws.Range("D20").Value = 64.1
$ws.Range("E20").Value = 63.9
$ws.Range("F20").Value = 63.8
$ws.Range("G20").Value = 63.8
$ws.Range("H20").Value = 64.2
$ws.Range("I20").Value = 65
$ws.Range("J20").Value = 66.1
$ws.Range("K20").Value = 67.6
$ws.Range("L20").Value = 69.3

# Row 21: West Pokot
$ws.Range("B21").Value = 54
$ws.Range("C21").Value = 52.9
$ws.Range("D21").Value = 51.8
$ws.Range("E21").Value = 50.8
$ws.Range("F21").Value = 49.8
$ws.Range("G21").Value = 48.8
$ws.Range("H21").Value = 48
$ws.Range("I21").Value = 47.3
$ws.Range("J21").Value = 46.8
$ws.Range("K21").Value = 46.6
$ws.Range("L21").Value = 46.7

